# Generate Report for Handoff
#
# The source file "d99b4360-...md" was renamed/regenerated as
# "5e1a45c3-...md" (new xliff hashes + a later handoff timestamp) and a
# second source file "ffff9b9a05c8-...md" became "Ready for handoff".
# Reflect both changes on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad4c595b28f4fe687a57a5b77a6a9c37bf0a65e4/e2e/"

$renamedFile = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md"
$newFile     = "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md"

$zhXlf = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.zh-cn.xlf"
$deXlf = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.de-de.xlf"

$zhHandoffDate = "2016-10-24 09:44:55"
$deHandoffDate = "2016-10-24 09:45:07"

# a lone "'" forces Excel to store an empty string in the cell instead of
# leaving it blank (COM drops a plain "" assignment entirely)
$blank = "'"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = $newFile
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = $blank
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = $deHandoffDate
$ws.Range("G3").NumberFormat = $dateFmt

$ws.Hyperlinks.Add($ws.Range("B3"), ($repoBase + $newFile), [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $newFile))

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# row 2 now reflects the renamed source file
$ws.Range("G2").Value = $zhXlf
$ws.Range("H2").Value = $zhHandoffDate
$ws.Range("H2").NumberFormat = $dateFmt

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), ($repoBase + $renamedFile), [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $renamedFile))

# row 3 is the brand new "Ready for handoff" entry
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = $zhXlf
$ws.Range("H3").Value = $zhHandoffDate
$ws.Range("H3").NumberFormat = $dateFmt
$ws.Range("I3").Value = $blank
$ws.Range("J3").Value = $blank
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = $dateFmt
$ws.Range("L3").Value = $blank
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = $blank
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = $blank

$ws.Hyperlinks.Add($ws.Range("A3"), ($repoBase + $newFile), [System.Type]::Missing, [System.Type]::Missing, $newFile)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# row 2 now reflects the renamed source file
$ws.Range("G2").Value = $deXlf
$ws.Range("H2").Value = $deHandoffDate
$ws.Range("H2").NumberFormat = $dateFmt

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), ($repoBase + $renamedFile), [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $renamedFile))

# row 3 is the brand new "Ready for handoff" entry
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = $deXlf
$ws.Range("H3").Value = $deHandoffDate
$ws.Range("H3").NumberFormat = $dateFmt
$ws.Range("I3").Value = $blank
$ws.Range("J3").Value = $blank
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = $dateFmt
$ws.Range("L3").Value = $blank
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = $blank
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = $blank

$ws.Hyperlinks.Add($ws.Range("A3"), ($repoBase + $newFile), [System.Type]::Missing, [System.Type]::Missing, $newFile)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))
